# Generate Report for Handback
# This script mirrors the "handback" report generation: it marks the
# zh-cn and de-de localization rows as handed back, fills in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns, adds matching hyperlinks and widens a few columns
# so the new content is readable.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderlineColor = 15570276   # OLE (BGR) form of RGB(0x64,0x95,0xED) == FF6495ED

function Set-HyperlinkLook($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $hyperlinkUnderlineColor
}

# ---------------------------------------------------------------
# Overview sheet: "Status" column text changes for both rows.
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

$zh.Range("I2").Value = "33acdf3f-5048-465a-acb9-fd83ae475c39.md"
$zh.Range("J2").Value = "33acdf3f-5048-465a-acb9-fd83ae475c39.d37d9c841bfde55415f5f4fc13fe3a0fc2d337ee.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-01 18:34:25"

$zh.Range("I3").Value = "dd09c3a5-3fc9-429e-a5ae-22ab3bb7ccaa.md"
$zh.Range("J3").Value = "dd09c3a5-3fc9-429e-a5ae-22ab3bb7ccaa.e95487c0cfcd4739ca7582b7feceb65529b014a6.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-01 18:34:25"

$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a5616c3186b6481ecb18b526f0e3aff308117cd/e2e/33acdf3f-5048-465a-acb9-fd83ae475c39.md", "", "", "33acdf3f-5048-465a-acb9-fd83ae475c39.md")
Set-HyperlinkLook $zh.Range("I2")

$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a5616c3186b6481ecb18b526f0e3aff308117cd/e2e/dd09c3a5-3fc9-429e-a5ae-22ab3bb7ccaa.md", "", "", "dd09c3a5-3fc9-429e-a5ae-22ab3bb7ccaa.md")
Set-HyperlinkLook $zh.Range("I3")

$zh.Columns.Item(3).ColumnWidth = 29.1666666666667
$zh.Columns.Item(9).ColumnWidth = 39.1666666666667
$zh.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

$de.Range("I2").Value = "33acdf3f-5048-465a-acb9-fd83ae475c39.md"
$de.Range("J2").Value = "33acdf3f-5048-465a-acb9-fd83ae475c39.d37d9c841bfde55415f5f4fc13fe3a0fc2d337ee.de-de.xlf"
$de.Range("K2").Value = "2016-09-01 18:34:33"

$de.Range("I3").Value = "dd09c3a5-3fc9-429e-a5ae-22ab3bb7ccaa.md"
$de.Range("J3").Value = "dd09c3a5-3fc9-429e-a5ae-22ab3bb7ccaa.e95487c0cfcd4739ca7582b7feceb65529b014a6.de-de.xlf"
$de.Range("K3").Value = "2016-09-01 18:34:33"

$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a5616c3186b6481ecb18b526f0e3aff308117cd/e2e/33acdf3f-5048-465a-acb9-fd83ae475c39.md", "", "", "33acdf3f-5048-465a-acb9-fd83ae475c39.md")
Set-HyperlinkLook $de.Range("I2")

$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a5616c3186b6481ecb18b526f0e3aff308117cd/e2e/dd09c3a5-3fc9-429e-a5ae-22ab3bb7ccaa.md", "", "", "dd09c3a5-3fc9-429e-a5ae-22ab3bb7ccaa.md")
Set-HyperlinkLook $de.Range("I3")

$de.Columns.Item(3).ColumnWidth = 29.1666666666667
$de.Columns.Item(9).ColumnWidth = 39.1666666666667
$de.Columns.Item(10).ColumnWidth = 39.1666666666667
